$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number (e.g. "215.11") need
# to be forced to Text first, otherwise Excel auto-converts them to a numeric
# value and loses the original formatting (trailing zeros, etc.).
$ws.Range("D2").Value = "27.585.58"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").Value = "1.664.38"
$ws.Range("E3").Value = "  -3.70%  "
$ws.Range("E4").Value = "  +0.15%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "215.11"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("E9").Value = "  -0.76%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0620"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("D12").Value = "1.899.77"
$ws.Range("E12").Value = "  -3.68%  "
$ws.Range("D13").Value = "1.647.38"
$ws.Range("E13").Value = "  -4.50%  "
$ws.Range("E14").Value = "  -3.36%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.560"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "27.579.86"
$ws.Range("E17").Value = "  -2.46%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "241.99"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").Value = "0.0$([char]0x2083)0730"
$ws.Range("E19").Value = "  -3.58%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.63"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -5.11%  "
$ws.Range("E21").Value = "  +0.21%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.48"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.48%  "
$ws.Range("E23").Value = "  -3.85%  "
$ws.Range("E24").Value = "  -2.79%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "147.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("E26").Value = "  -4.31%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "16.40"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("E32").Value = "  -2.83%  "
$ws.Range("D33").Value = "1.458.69"
$ws.Range("E33").Value = "  -3.05%  "
$ws.Range("E34").Value = "  -4.97%  "
$ws.Range("E35").Value = "  -5.06%  "
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("E37").Value = "  -4.63%  "
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("E39").Value = "  -5.64%  "
$ws.Range("E40").Value = "  -2.85%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "69.52"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.11%  "
$ws.Range("E42").Value = "  +0.17%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.42"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -6.03%  "
$ws.Range("E44").Value = "  -3.44%  "
$ws.Range("D45").Value = "1.807.56"
$ws.Range("E45").Value = "  -3.71%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.789"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("E47").Value = "  -1.23%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "89.10"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("E49").Value = "  -3.51%  "
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("E51").Value = "  -4.82%  "
